# Refresh crypto price/volume data (and fix two mis-ordered rows) as
# captured by the scheduled "Updated cryptos list" GitHub Actions run.
#
# Every cell written here is a plain text cell in the workbook (column D
# holds price strings like "2.414.07" and column E holds padded percent
# strings like "  +0.95%  "). Excel's COM layer auto-detects plausible
# numbers/percents typed into a cell and silently coerces them, so each
# write is given a leading apostrophe to force literal text, then the
# cell style is reset to "Normal" so no stray quote-prefix formatting is
# left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'57.118.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.95%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.414.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.97%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'488.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.44%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'154.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.17%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.34%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'0.605"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +19.26%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'2.431.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.30%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'6.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +10.98%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.100"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.95%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.331"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.33%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("E13").Value = "'  +1.17%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'2.827.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.26%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'57.159.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.81%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'20.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.97%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.99%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'2.430.44"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.60%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'4.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.68%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'320.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.54%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'10.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.01%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("E22").Value = "'  -0.20%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'5.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.99%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'58.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.40%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("B25").Value = "'Polygon"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'0.402"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.38%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("B26").Value = "'Binance-PegBSC-USD"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value = "'0.990"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.47%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.160"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.37%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'2.524.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.22%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'7.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.97%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'0.0₃0789"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.69%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.09%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'150.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.39%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'18.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +3.83%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("E34").Value = "'  +0.94%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'5.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.36%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'3.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.45%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'1.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.28%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("D38").Value = "'0.817"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -7.35%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'34.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.10%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'1.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.05%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'3.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.84%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'  +5.09%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("B43").Value = "'Bittensor"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'277.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +5.05%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("B44").Value = "'FirstDigitalUSD"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.994"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.47%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.595"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.37%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.0535"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.73%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "'  -0.19%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.0229"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.73%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'4.61"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.30%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'17.87"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.48%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("B51").Value = "'SuiNetwork"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.678"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +6.29%  "
$ws.Range("E51").Style = "Normal"

